$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Alterei alguns use cases: o ator do caso de uso "Carro Pronto" passa de
# "Administrador" para "Funcionário".
$ws.Range("C3").Value = "Funcionário"
